$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry describes the cell values that changed for a given row,
# in columns B (Coin), C (Link), D (Price) and E (Volume(1h)).
$updates = @(
    @{Row=2; D='51.924.25'; E='  +0.16%  '},
    @{Row=3; D='2.780.59'; E='  -1.16%  '},
    @{Row=4; E='  +0.03%  '},
    @{Row=5; D='357.49'; E='  +1.00%  '},
    @{Row=6; D='109.35'; E='  -3.47%  '},
    @{Row=7; D='0.565'; E='  +2.33%  '},
    @{Row=8; D='1.00'; E='  +0.05%  '},
    @{Row=9; E='  -0.80%  '},
    @{Row=10; D='40.10'; E='  -3.52%  '},
    @{Row=11; E='  +0.15%  '},
    @{Row=12; E='  +0.73%  '},
    @{Row=13; D='19.41'; E='  -2.77%  '},
    @{Row=14; D='7.62'; E='  -1.01%  '},
    @{Row=15; D='3.218.47'; E='  -0.47%  '},
    @{Row=16; D='2.880.07'; E='  +2.36%  '},
    @{Row=17; D='0.928'; E='  +3.68%  '},
    @{Row=18; D='51.873.90'; E='  +0.24%  '},
    @{Row=19; E='  +0.58%  '},
    @{Row=20; D='3.16'; E='  +0.34%  '},
    @{Row=21; D='13.05'; E='  -3.36%  '},
    @{Row=22; E='  -1.51%  '},
    @{Row=23; D='274.28'; E='  +1.59%  '},
    @{Row=24; D='69.96'; E='  +0.25%  '},
    @{Row=25; D='2.73'; E='  -2.03%  '},
    @{Row=26; D='26.55'; E='  -0.51%  '},
    @{Row=27; E='  -0.02%  '},
    @{Row=28; D='10.15'; E='  -1.31%  '},
    @{Row=29; E='  -1.20%  '},
    @{Row=30; E='  +2.26%  '},
    @{Row=31; B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='51.66'; E='  +2.07%  '},
    @{Row=32; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.0464'; E='  +2.95%  '},
    @{Row=33; D='34.00'},
    @{Row=34; D='5.70'; E='  -1.97%  '},
    @{Row=35; E='  +1.68%  '},
    @{Row=36; D='5.25'; E='  +7.29%  '},
    @{Row=37; E='  +0.03%  '},
    @{Row=38; E='  +0.91%  '},
    @{Row=39; D='18.11'; E='  -0.89%  '},
    @{Row=40; E='  -4.29%  '},
    @{Row=41; E='  -1.04%  '},
    @{Row=42; E='  -0.39%  '},
    @{Row=43; E='  -3.04%  '},
    @{Row=44; D='121.25'; E='  -4.18%  '},
    @{Row=45; D='22.00'; E='  -6.93%  '},
    @{Row=46; D='2.070.46'},
    @{Row=47; E='  -2.89%  '},
    @{Row=48; D='2.24'; E='  -2.95%  '},
    @{Row=49; D='5.69'; E='  +0.46%  '},
    @{Row=50; D='0.919'; E='  -2.28%  '},
    @{Row=51; D='8.95'; E='  +0.54%  '}
)

foreach ($u in $updates) {
    $row = $u.Row
    foreach ($col in 'B','C','D','E') {
        if ($u.ContainsKey($col)) {
            $cell = $ws.Range("$col$row")
            if ($col -eq 'D') {
                # Price values such as "1.00", "51.66" or "34.00" look like
                # numbers/dates to Excel's auto-detection. Force the cell to
                # Text format while assigning so the literal string is kept,
                # then restore the default "Normal" style so no stray
                # number-format style is left behind on the cell.
                $cell.NumberFormat = "@"
                $cell.Value = $u[$col]
                $cell.Style = "Normal"
            } else {
                $cell.Value = $u[$col]
            }
        }
    }
}
